# means_t-results.xlsx adjustments:
# - exclude South Africa and Netherlands rows (4th/5th row of every 5-row group)
# - include omega -> updated "group overall" value (row 1 of each group) recalculated
# - switch to one-sided testing -> updated SDs / ancillary values
# - sheet view / column width cosmetic resets

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Delete the 4th and 5th data row of each original 5-row group
#    (original rows: 2-6, 7-11, 12-16, ..., 47-51). Delete bottom-up so
#    earlier row numbers stay valid while iterating.
$rowsToDelete = @(5,6,10,11,15,16,20,21,25,26,30,31,35,36,40,41,45,46,50,51)
$sorted = $rowsToDelete | Sort-Object -Descending
foreach ($r in $sorted) {
    $ws.Rows($r).Delete()
}

# 2) Update the recomputed summary value (column B) for the first row of
#    each remaining 3-row group.
$ws.Range("B2").Value  = "3.91 (0.75)"
$ws.Range("B5").Value  = "4 (0.78)"
$ws.Range("B8").Value  = "3.98 (0.76)"
$ws.Range("B11").Value = "4 (0.71)"
$ws.Range("B14").Value = "35.33 (25.76)"
$ws.Range("B17").Value = "36.28 (26.32)"
$ws.Range("B20").Value = "4.11 (2.03)"
$ws.Range("B23").Value = "4.12 (1.98)"
$ws.Range("B26").Value = "3.97 (2.03)"
$ws.Range("B29").Value = "4.07 (1.99)"

# 3) Column widths: column A narrower, column B back to the sheet default
#    (no custom width).
$ws.Columns("A").ColumnWidth = 18.6328125
$ws.Columns("B").EntireColumn.AutoFit() | Out-Null
$ws.Columns("B").ColumnWidth = 8.7265625

# 4) Selection / view: land on A7 with nothing scrolled away.
$ws.Range("A7").Select()
